# Apply updated crypto price/volume data to Sheet1 (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.135.19'
$ws.Range("E2").Value = '  +0.61%  '

$ws.Range("D3").Value = '1.565.92'
$ws.Range("E3").Value = '  +0.52%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.94%  '

$ws.Range("E6").Value = '  +0.41%  '

$ws.Range("E7").Value = '  +0.68%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.94'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.249'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.19%  '

$ws.Range("E10").Value = '  +0.46%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0865'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.02%  '

$ws.Range("D12").Value = '1.789.33'
$ws.Range("E12").Value = '  +0.46%  '

$ws.Range("D13").Value = '1.540.00'
$ws.Range("E13").Value = '  +0.38%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.78'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.37%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.518'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.56%  '

$ws.Range("D16").Value = '27.107.29'
$ws.Range("E16").Value = '  +0.52%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.13'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.29%  '

$ws.Range("D18").Value = '0.0₃0701'
$ws.Range("E18").Value = '  -0.86%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.16'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.69%  '

$ws.Range("E21").Value = '  +0.71%  '

$ws.Range("E22").Value = '  +0.93%  '

$ws.Range("E23").Value = '  -0.64%  '

$ws.Range("E24").Value = '  +0.84%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.16'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.61'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.02%  '

$ws.Range("E28").Value = '  +1.72%  '

$ws.Range("E29").Value = '  +0.70%  '

$ws.Range("E30").Value = '  +4.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0471'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.42%  '

$ws.Range("E32").Value = '  +0.12%  '

$ws.Range("E33").Value = '  +2.10%  '

$ws.Range("D34").Value = '1.449.37'
$ws.Range("E34").Value = '  +1.74%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.09'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.92%  '

$ws.Range("E36").Value = '  -0.07%  '

$ws.Range("E37").Value = '  +1.87%  '

$ws.Range("E38").Value = '  +0.96%  '

$ws.Range("E39").Value = '  +0.22%  '

$ws.Range("E40").Value = '  +2.20%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.806'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.23%  '

$ws.Range("E42").Value = '  +0.67%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.34'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.82%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.50'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.57%  '

$ws.Range("E46").Value = '  -0.40%  '

$ws.Range("D47").Value = '1.701.72'
$ws.Range("E47").Value = '  +0.50%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.91'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.97%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0518'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.44%  '

$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₆0100'
$ws.Range("E50").Value = '  +1.69%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0958'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.16%  '

